# Refresh the crypto price/volume table with the latest scrape (6 Jan 2023,
# hour 8) and fix the BKEXToken/KickToken row ordering that had been swapped.
#
# All figures in column D (Price) and column E (Volume(1h)) are stored as
# plain TEXT in this sheet (not numbers), e.g. "257.06" and "0.07%". Writing
# a numeric- or percent-looking string straight into .Value would make Excel
# silently reinterpret it as a real number/percentage, which would corrupt
# both the stored value and the cell's number format. To avoid that we force
# the cell to Text format before writing, then drop the style back to
# "Normal" afterwards so we don't leave a stray number-format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Coin name / link swap (rows were in the wrong order)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

# row, price (D), volume% (E) -- $null means "leave unchanged"
$rows = @(
    @{ R = 2;  D = "256.91";     E = "-0.03%" }
    @{ R = 3;  D = "27.02";      E = "-0.49%" }
    @{ R = 4;  D = "4.550";      E = "-5.08%" }
    @{ R = 5;  D = "0.05897";    E = $null }
    @{ R = 6;  D = "6.609";      E = "-0.74%" }
    @{ R = 7;  D = "0.8499";     E = "-2.34%" }
    @{ R = 8;  D = "0.9312";     E = "-2.09%" }
    @{ R = 9;  D = "0.1375";     E = "-2.06%" }
    @{ R = 10; D = "0.04240";    E = "11.13%" }
    @{ R = 11; D = "0.07018";    E = "-1.98%" }
    @{ R = 12; D = "0.03045";    E = "-4.79%" }
    @{ R = 13; D = "0.09105";    E = "-1.68%" }
    @{ R = 14; D = "0.001526";   E = "-1.16%" }
    @{ R = 15; D = "0.0006037";  E = "-94.24%" }
    @{ R = 16; D = "0.006083";   E = "1.20%" }
    @{ R = 17; D = $null;        E = "-0.36%" }
    @{ R = 18; D = "3.172";      E = "-0.66%" }
    @{ R = 19; D = $null;        E = "-1.23%" }
    @{ R = 20; D = $null;        E = "-1.75%" }
    @{ R = 21; D = $null;        E = "-0.97%" }
    @{ R = 22; D = "3.906";      E = "2.14%" }
    @{ R = 23; D = $null;        E = "1.39%" }
    @{ R = 24; D = "0.001221";   E = "-0.26%" }
    @{ R = 25; D = $null;        E = "-4.39%" }
    @{ R = 26; D = $null;        E = "-0.05%" }
    @{ R = 27; D = $null;        E = "1.99%" }
    @{ R = 40; D = "0.03793";    E = "-1.24%" }
    @{ R = 41; D = "0.006297";   E = "0.75%" }
    @{ R = 42; D = "0.1099";     E = "-0.15%" }
    @{ R = 43; D = "0.002199";   E = "-2.35%" }
    @{ R = 44; D = $null;        E = "32.45%" }
    @{ R = 45; D = "0.00005354"; E = "-2.65%" }
    @{ R = 46; D = $null;        E = "-0.05%" }
    @{ R = 47; D = "0.05097";    E = "-42.41%" }
    @{ R = 48; D = "0.2521";     E = "10,478.70%" }
    @{ R = 49; D = $null;        E = "-0.05%" }
    @{ R = 50; D = $null;        E = "-0.05%" }
)

foreach ($row in $rows) {
    if ($null -ne $row.D) { Set-TextValue "D$($row.R)" $row.D }
    if ($null -ne $row.E) { Set-TextValue "E$($row.R)" $row.E }
}
